$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 2 new rows into the data table (before the old last row,
#    which was row 19) so the table grows from 4 rows to 6 rows.
#    This naturally pushes the old row 19 down to row 21, and the
#    blank spacer rows + signature footer (old rows 24-25) down to
#    rows 26-27.
# ------------------------------------------------------------------
$ws.Range("A19:A20").EntireRow.Insert(-4121)   # xlShiftDown

# Seed the two new rows (19 and 20) with the formatting/content of an
# existing "normal" data row (row 18) so they pick up the same cell
# styles/borders as the rest of the table instead of plain defaults.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("B18:J18").Copy($ws.Range("B20:J20"))

# ------------------------------------------------------------------
# 2. Rewrite the six data rows (16-21) with the updated dataset: the
#    two workers (CC 45592082 / CC 22810295) now each have three
#    "Periodo Mora" entries (2506, 2507 and the new 2508), grouped by
#    period.
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45592082"
$ws.Range("D16").Value = "MARINELA BLANQUICETT BARRERA"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 58476
$ws.Range("G16").Value = 1461908

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "22810295"
$ws.Range("D17").Value = "ZULAY DEL CARMEN ROMERO SANCHEZ"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45592082"
$ws.Range("D18").Value = "MARINELA BLANQUICETT BARRERA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 58476
$ws.Range("G18").Value = 1461908

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "22810295"
$ws.Range("D19").Value = "ZULAY DEL CARMEN ROMERO SANCHEZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45592082"
$ws.Range("D20").Value = "MARINELA BLANQUICETT BARRERA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 58476
$ws.Range("G20").Value = 1461908

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "22810295"
$ws.Range("D21").Value = "ZULAY DEL CARMEN ROMERO SANCHEZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# ------------------------------------------------------------------
# 3. Update the summary fields at the top of the sheet: total "Valor
#    Mora" (now the sum over 6 rows instead of 4) and the updated
#    "Cant. Periodos" count (3 distinct periods instead of 2).
# ------------------------------------------------------------------
$ws.Range("E11").Value = 346248
$ws.Range("F13").Value = 3
